$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new results row (row 3)
$ws.Range("A3").Value = "vanilla"
$ws.Range("B3").Value = "2x200 F G H"
$ws.Range("C3").Value = 0.001
$ws.Range("D3").Value = 30
$ws.Range("E3").Value = 0.029861111111111113
$ws.Range("F3").Value = 0.85409500000000005
$ws.Range("G3").Value = 0.85148599999999997
$ws.Range("H3").Value = 0.85470400000000002
$ws.Range("I3").Value = 0.85199499999999995

# Update the active selection to I3
$ws.Range("I3").Select()
